$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need to be forced to
# text (format "@") so Excel does not silently convert them to a Double
# (which would lose trailing zeros / introduce float noise). The number
# format + style are reset back to "Normal" immediately after the write
# so the cell keeps its original (default) style index.
$numericLooking = @("D4", "D5", "D7", "D8", "D9", "D11", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D22", "D25", "D26", "D27", "D29", "D31", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $numericLooking) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.644.18"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "1.802.09"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "315.66"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").Value = "0.5321"
$ws.Range("E7").Value = "  -7.32%  "
$ws.Range("D8").Value = "0.3772"
$ws.Range("E8").Value = "  -1.96%  "
$ws.Range("D9").Value = "42.49"
$ws.Range("E9").Value = "  -2.13%  "
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("D11").Value = "1.112"
$ws.Range("E11").Value = "  -2.29%  "
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "20.69"
$ws.Range("E13").Value = "  -3.05%  "
$ws.Range("D14").Value = "6.151"
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("D15").Value = "7.348"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "1.796.76"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "90.20"
$ws.Range("E17").Value = "  -2.29%  "
$ws.Range("D18").Value = "0.00001064"
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("D19").Value = "0.06459"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("D22").Value = "5.897"
$ws.Range("E22").Value = "  -1.73%  "
$ws.Range("D23").Value = "28.556.03"
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("E24").Value = "  -2.87%  "
$ws.Range("D25").Value = "2.095"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "159.75"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").Value = "20.41"
$ws.Range("E27").Value = "  -2.57%  "
$ws.Range("D28").Value = "2.006.84"
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").Value = "2.352"
$ws.Range("E29").Value = "  -2.43%  "
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("D31").Value = "1.099"
$ws.Range("E31").Value = "  -5.24%  "
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("D33").Value = "3.700"
$ws.Range("E33").Value = "  +2.13%  "
$ws.Range("D34").Value = "5.637"
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("E35").Value = "  +4.07%  "
$ws.Range("D36").Value = "0.06413"
$ws.Range("E36").Value = "  +5.32%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.02303"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "8.789"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").Value = "5.033"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").Value = "1.217"
$ws.Range("E40").Value = "  +4.83%  "
$ws.Range("D41").Value = "11.27"
$ws.Range("E41").Value = "  -3.87%  "
$ws.Range("D42").Value = "0.6211"
$ws.Range("E42").Value = "  -3.56%  "
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").Value = "1.407"
$ws.Range("E44").Value = "  +2.11%  "
$ws.Range("D45").Value = "13.32"
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.5848"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.687"
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("D48").Value = "125.71"
$ws.Range("E48").Value = "  +2.85%  "
$ws.Range("D49").Value = "1.940"
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("D50").Value = "1.148"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").Value = "0.06885"
$ws.Range("E51").Value = "  +0.42%  "

foreach ($ref in $numericLooking) {
    $ws.Range($ref).Style = "Normal"
}
